$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 81
$ws1.Range("F8").Value = 123
$ws1.Range("F9").Value = 8941
$ws1.Range("F11").Value = 333
$ws1.Range("F17").Value = 242
$ws1.Range("F21").Value = 1111

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 3
$ws2.Range("F3").Value = 5
$ws2.Range("G3").Value = 138

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3
$ws4.Range("F6").Value = 81
$ws4.Range("F8").Value = 5
$ws4.Range("G8").Value = 138
$ws4.Range("F10").Value = 123
$ws4.Range("F11").Value = 8941
$ws4.Range("F13").Value = 333
$ws4.Range("F19").Value = 242
$ws4.Range("F23").Value = 1111
